# 17th Commit - Scanning Operational
# Zero-out the stale per-column scan values on the Zmatrix sheet (row 1-3)
# now that the Ras-Pi/Arducam scanning pipeline is live but not yet
# validated against the laser line.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Zmatrix")

$cellsToClear = @(
    "H1","I1","J1","K1","L1","M1","N1","O1",
    "AM1","AN1","AO1","AQ1","AT1","AU1",
    "J2","K2","L2","M2","O2","P2",
    "K3","L3","N3","O3","X3","AK3",
    "AN3","AO3","AP3","AQ3","AR3","AS3","AT3","AU3","AV3"
)

foreach ($addr in $cellsToClear) {
    $ws.Range($addr).Value = 0
}
